$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '246.20'
Set-TextValue 'D3' '24.16'
Set-TextValue 'D4' '5.274'
Set-TextValue 'D5' '0.05803'
Set-TextValue 'D6' '6.499'
Set-TextValue 'D7' '3.131'
Set-TextValue 'D8' '0.8175'
Set-TextValue 'D9' '0.8539'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue 'D10' '0.1360'
$ws.Range("E10").Value = '9WazirXWRX'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue 'D11' '0.06929'
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue 'D12' '0.03130'
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue 'D13' '0.02868'
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue 'D14' '0.09398'
$ws.Range("E14").Value = '13BitMartTokenBMX'
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue 'D15' '3.739'
$ws.Range("E15").Value = '14MCDexMCB'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue 'D16' '0.001518'
$ws.Range("E16").Value = '15BitForexTokenBF'
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue 'D17' '0.04680'
$ws.Range("E17").Value = '16CoinExTokenCET'
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue 'D18' '0.0005982'
$ws.Range("E18").Value = '17OneONE'
Set-TextValue 'D19' '0.006264'
Set-TextValue 'D20' '0.001233'
Set-TextValue 'D21' '0.004629'
Set-TextValue 'D23' '3.504'
Set-TextValue 'D25' '0.3191'
Set-TextValue 'D26' '0.1345'
Set-TextValue 'D27' '0.1357'
Set-TextValue 'D28' '0.0002328'
Set-TextValue 'D40' '0.03667'
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue 'D41' '0.1058'
$ws.Range("E41").Value = '40BKEXTokenBKK'
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue 'D42' '0.002748'
$ws.Range("E42").Value = '41CEJICEJIBestin24h'
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue 'D43' '0.003014'
$ws.Range("E43").Value = '42KickTokenKICKWorstin24h'
Set-TextValue 'D44' '0.007475'
Set-TextValue 'D45' '0.00005266'
Set-TextValue 'D47' '0.3695'
$ws.Range("E47").Value = '46CoinbaseStockTokenCOIN'
Set-TextValue 'D48' '0.002244'
Set-TextValue 'D49' '0.00002097'
Set-TextValue 'D50' '0.0001997'
